# NSE Daily update — appends a new day's block of 5 stock rows
# (RELIANCE.NS, TCS.NS, HDFCBANK.NS, BHARTIARTL.NS, ICICIBANK.NS) to the
# "NSE_2024-12-08" sheet, duplicating the most-recently fetched block
# (rows 7-11) into new rows 12-16, the way the daily scraper job appends
# a fresh fetch each day. One field (RELIANCE's lastDividendValue, column
# BC) changed between fetches, so it is overridden after the copy.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NSE_2024-12-08")

$lastCol = 87   # column CI
$rowPairs = @(
    @{ Src = 7;  Dst = 12 },
    @{ Src = 8;  Dst = 13 },
    @{ Src = 9;  Dst = 14 },
    @{ Src = 10; Dst = 15 },
    @{ Src = 11; Dst = 16 }
)

foreach ($pair in $rowPairs) {
    $srcRow = $pair.Src
    $dstRow = $pair.Dst

    for ($col = 1; $col -le $lastCol; $col++) {
        $srcCell = $ws.Cells.Item($srcRow, $col)
        $dstCell = $ws.Cells.Item($dstRow, $col)

        $v = $srcCell.Value2

        if ($null -eq $v) {
            # Source column wasn't populated for this stock either - leave blank.
            continue
        }

        $typeName = $v.GetType().Name

        if ($typeName -eq "String") {
            # Force literal text so date-like strings ("2024-12-07") and
            # blank placeholders aren't reinterpreted as dates/numbers.
            if ($v -eq "") {
                $dstCell.Value2 = "'"
            } else {
                $dstCell.Value2 = "'" + $v
            }
        } else {
            $dstCell.Value2 = $v
        }
    }
}

# RELIANCE.NS lastDividendValue refreshed from 10 to 5 in this fetch.
$ws.Cells.Item(12, 55).Value2 = 5
